$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 258, shifting existing rows 258:270 down to 259:271.
$ws.Rows.Item(258).Insert()

# Populate the new row with the new weekly record.
$ws.Cells.Item(258, 1).Value = 4
$ws.Cells.Item(258, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(258, 3).Value = "Los Lagos"
$ws.Cells.Item(258, 4).Value = 44746
$ws.Cells.Item(258, 5).Value = 10
$ws.Cells.Item(258, 6).Value = 100112017
$ws.Cells.Item(258, 7).Value = "Apio"
$ws.Cells.Item(258, 8).Value = "Americana (o)"
$ws.Cells.Item(258, 9).Value = "Primera"
$ws.Cells.Item(258, 10).Value = 25
$ws.Cells.Item(258, 11).Value = 12000
$ws.Cells.Item(258, 12).Value = 12000
$ws.Cells.Item(258, 13).Value = 12000
$ws.Cells.Item(258, 14).Value = "`$/docena de matas"
$ws.Cells.Item(258, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(258, 16).Value = 2000
$ws.Cells.Item(258, 17).Value = 6
$ws.Cells.Item(258, 18).Value = "Hortaliza"
